$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table shrinks from a 4x4 grid (Source/A/FFR/LF headers, A Lag/FFR
# Lag/LF Lag rows) to a 3x3 grid (Source/FFR/LF headers, FFR Lag/LF Lag
# rows) with refreshed coefficient values. Clear the old range first so
# column D / row 4 don't linger, then write the new table.

$ws.Range("A1:D4").ClearContents()

$ws.Range("A1").Value = "Source"
$ws.Range("B1").Value = "FFR"
$ws.Range("C1").Value = "LF"

$ws.Range("A2").Value = "FFR Lag"
$ws.Range("A3").Value = "LF Lag"

$ws.Range("B2").Value = "1.88***"
$ws.Range("B3").Value = "3.77*"

$ws.Range("C2").Value = "0.47***"
$ws.Range("C3").Value = "0.76*"

$ws.Columns("D").Delete()
$ws.Rows(4).Delete()
